$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for columns G and H ---
$ws.Range("G1").Value = "Ny leverandør"
$ws.Range("H1").Value = "TCV_range"

# Column F header changes meaning from "TCV_range" to "Årsag"
$ws.Range("F1").Value = "Årsag"

# Copy the bold/centered header formatting from F1 onto the new G1/H1 headers
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Per-row data for the new "Årsag" (F) and "Ny leverandør" (G) columns ---
$data = @(
    @{ Row=2; F='Pris'; G='' }
    @{ Row=3; F='Bruger ikke produktet'; G='' }
    @{ Row=4; F='Anden årsag (angiv hvilken i bemærkninger)'; G='Danløn' }
    @{ Row=5; F='Anden årsag (angiv hvilken i bemærkninger)'; G='Lessor' }
    @{ Row=6; F='Strategisk beslutning'; G='' }
    @{ Row=7; F='Bruger ikke produktet'; G='' }
    @{ Row=8; F='Fusionerer med anden virksomhed'; G='' }
    @{ Row=9; F='Bruger ikke produktet'; G='' }
    @{ Row=10; F='Fusionerer med anden virksomhed'; G='' }
    @{ Row=11; F='Ikke oplyst'; G='' }
    @{ Row=12; F='Virksomheden lukker'; G='' }
    @{ Row=13; F='Bruger ikke produktet'; G='' }
    @{ Row=14; F='Ikke oplyst'; G='' }
    @{ Row=15; F='Utilfredshed (Service - uddyb i bemærkninger)'; G='' }
    @{ Row=16; F='Ikke oplyst'; G='' }
    @{ Row=17; F='Ikke oplyst'; G='' }
    @{ Row=18; F='Ikke oplyst'; G='' }
    @{ Row=19; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=20; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=21; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=22; F='Bruger ikke produktet'; G='DataLøn' }
    @{ Row=23; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=24; F='Strategisk beslutning'; G='' }
    @{ Row=25; F='Ikke oplyst'; G='' }
    @{ Row=26; F='Virksomheden lukker'; G='' }
    @{ Row=27; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=28; F='Anden årsag (angiv hvilken i bemærkninger)'; G='' }
    @{ Row=29; F='Strategisk beslutning'; G='' }
    @{ Row=30; F='Virksomheden lukker'; G='' }
    @{ Row=31; F='Strategisk beslutning'; G='' }
    @{ Row=32; F='Utilfredshed (Service - uddyb i bemærkninger)'; G='' }
    @{ Row=33; F='Ikke oplyst'; G='' }
    @{ Row=34; F='Ikke oplyst'; G='' }
    @{ Row=35; F='Ikke oplyst'; G='' }
    @{ Row=36; F='Bruger ikke produktet'; G='' }
    @{ Row=37; F='Bruger ikke produktet'; G='' }
    @{ Row=38; F='Bruger ikke produktet'; G='' }
    @{ Row=39; F='Ikke oplyst'; G='' }
)

foreach ($d in $data) {
    $row = $d.Row
    $ws.Cells.Item($row, 6).Value = $d.F
    if ($d.G -ne "") {
        $ws.Cells.Item($row, 7).Value = $d.G
    }
    $ws.Cells.Item($row, 8).Value = "0-20000"
}
